$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.ClearFormats()
}

Set-TextValue 2 4 '288.54'
Set-TextValue 2 5 '1.20%'
Set-TextValue 3 4 '29.20'
Set-TextValue 3 5 '1.99%'
Set-TextValue 4 4 '5.079'
Set-TextValue 4 5 '3.07%'
Set-TextValue 5 4 '0.06678'
Set-TextValue 5 5 '2.79%'
Set-TextValue 6 4 '7.355'
Set-TextValue 6 5 '1.97%'
Set-TextValue 7 2 'GateToken'
Set-TextValue 7 3 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 7 4 '3.406'
Set-TextValue 7 5 '1.06%'
Set-TextValue 8 2 'FTXToken'
Set-TextValue 8 3 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue 8 4 '1.349'
Set-TextValue 8 5 '1.71%'
Set-TextValue 9 2 'MXToken'
Set-TextValue 9 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 9 4 '0.9182'
Set-TextValue 9 5 '0.77%'
Set-TextValue 10 2 'WazirX'
Set-TextValue 10 3 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 10 4 '0.1579'
Set-TextValue 10 5 '2.72%'
Set-TextValue 11 2 'LiechtensteinCryptoassetsExchange'
Set-TextValue 11 3 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 11 4 '0.06805'
Set-TextValue 11 5 '5.71%'
Set-TextValue 12 2 'MandalaExchangeToken'
Set-TextValue 12 3 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 12 4 '0.07701'
Set-TextValue 12 5 '1.19%'
Set-TextValue 13 2 'BitrueCoin'
Set-TextValue 13 3 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 13 4 '0.02935'
Set-TextValue 13 5 '-1.74%'
Set-TextValue 14 2 'BitMartToken'
Set-TextValue 14 3 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 14 4 '0.08996'
Set-TextValue 14 5 '0.19%'
Set-TextValue 15 2 'BitForexToken'
Set-TextValue 15 3 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 15 4 '0.001571'
Set-TextValue 15 5 '-1.58%'
Set-TextValue 16 2 'CoinExToken'
Set-TextValue 16 3 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue 16 4 '0.04509'
Set-TextValue 16 5 '0.76%'
Set-TextValue 17 2 'One'
Set-TextValue 17 3 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 17 4 '0.0006452'
Set-TextValue 17 5 '-1.20%'
Set-TextValue 18 2 'TigerCash'
Set-TextValue 18 3 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 18 4 '0.006257'
Set-TextValue 18 5 '3.21%'
Set-TextValue 19 2 'LEO'
Set-TextValue 19 3 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 19 4 '3.454'
Set-TextValue 19 5 '-0.20%'
Set-TextValue 20 5 '-0.99%'
Set-TextValue 21 5 '2.03%'
Set-TextValue 22 4 '0.1309'
Set-TextValue 22 5 '-3.00%'
Set-TextValue 23 4 '4.071'
Set-TextValue 23 5 '2.50%'
Set-TextValue 24 4 '0.1582'
Set-TextValue 24 5 '1.68%'
Set-TextValue 25 4 '0.001191'
Set-TextValue 25 5 '0.24%'
Set-TextValue 26 4 '0.004116'
Set-TextValue 26 5 '-4.77%'
Set-TextValue 27 4 '0.0001199'
Set-TextValue 27 5 '1.50%'
Set-TextValue 28 4 '0.0001617'
Set-TextValue 28 5 '-1.25%'
Set-TextValue 40 4 '0.04202'
Set-TextValue 40 5 '1.14%'
Set-TextValue 41 4 '0.006707'
Set-TextValue 41 5 '0.04%'
Set-TextValue 42 4 '0.1239'
Set-TextValue 42 5 '0.56%'
Set-TextValue 43 4 '0.002090'
Set-TextValue 43 5 '-3.82%'
Set-TextValue 44 4 '0.01331'
Set-TextValue 44 5 '13.05%'
Set-TextValue 45 4 '0.00005721'
Set-TextValue 45 5 '5.14%'
Set-TextValue 46 4 '1.968'
Set-TextValue 46 5 '26.03%'
Set-TextValue 47 5 '-29.44%'
